# hs-logger_priorities.xlsx -- "473 Additional commands added."
#
# The "Points files" todo (row 37) got its Importance/Difficulty/Time-in-
# minutes numbers corrected, and -- as a side effect of Excel re-filling the
# Priority column's formula (column E) while the user was working in that
# area -- the three drag-fill chunks of E2*C2/D2-style formulas collapsed
# into shared-formula groups on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 ("Points files"): Importance/Difficulty/Time in minutes -------
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 5
$ws.Range("D37").Value = 2

# --- Re-enter the Priority formula over each fill-handle chunk -------------
# Re-assigning Formula across a whole range is what makes Excel collapse the
# per-row formulas into <f t="shared" .../> groups instead of one <f> per
# cell.
$ws.Range("E2:E33").Formula = "=B2*C2/D2"
$ws.Range("E34:E65").Formula = "=B34*C34/D34"
$ws.Range("E66:E97").Formula = "=B66*C66/D66"

# The last drag overshot the populated data (rows 66-76) down to row 97.
# Remove the phantom trailing rows again -- this leaves the shared formula's
# recorded ref="E66:E97" span in place on E66 (as in the saved file) while
# the sheet's actual used range/dimension goes back to row 76.
$ws.Range("E77:E97").ClearContents()
$ws.Range("E77:E97").Delete(-4162)

# --- Selection left where the user's cursor ended up after the edit -------
$ws.Range("A77").Select()
